$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting from the last existing data row (row 11) down to the three new rows
$ws.Range("A11:B11").Copy()
$ws.Range("A12:B12").PasteSpecial(-4122)
$ws.Range("A13:B13").PasteSpecial(-4122)
$ws.Range("A14:B14").PasteSpecial(-4122)

# Enter sequence values first (column B), then names (column A), matching
# the order new shared strings were authored in the source workbook.
$ws.Range("B12").Value = "CCCAACAATCCAACTCACGTTGAGATCCGCCAAAACGTCCACATTGGCTGTGTTCAGGCGGAGAGTCACCAACGTGGGACCACCGTCGTCGGTCTACACCGCCACCGTCCGAGCACCGAAAGGAGTAGAAATCACGGTGGAGCCACAGAGTTTGTCATTTTCAAAGGCTTCACAAAAGAGAAGCTTCAAAGTGGTGGT"
$ws.Range("B13").Value = "GAGCTGGAACATGTTCTGGTTGCAGCCGGAGGATCCTTGCCGCGTTTACAATCTCTGTGGTCAATTAGGGTTTTGTAGCAGCGAATTGCTCAAGCCCTGTGC"
$ws.Range("B14").Value = "GTATGCACGACTGGATCACTGAGAACCTCCGTGCGTGTGGCGGCACTTATCAGACATGTATCTGCGCCGTACCTTTCTTGGCAAAAAAGCAAGGTCTCGTGACCGTCACGTGCGATCCCAAGAACATCGAACACATGCTCAAGACCAGGTTCGACAACTACC"

$ws.Range("A12").Value = "XSP1"
$ws.Range("A13").Value = "RLK4"
$ws.Range("A14").Value = "CYP86A2"

# Restore active-cell selection like the original author left it
$ws.Range("D18").Select() | Out-Null
